# Adds a function for entering a single row of "Name / Address / Contact"
# data into the active sheet, then uses it to populate two more rows
# (row 6 fully, row 7 skipping the Address column) below the existing data.

function Add-SingleSheetDataEntry($Worksheet, $Row, $Name, $Address, $Contact) {
    if ($Name) {
        $Worksheet.Cells.Item($Row, 1).Value = $Name
    }
    if ($Address) {
        $Worksheet.Cells.Item($Row, 2).Value = $Address
    }
    if ($Contact) {
        $Worksheet.Cells.Item($Row, 3).Value = $Contact
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: full entry (Name, Address, Contact)
Add-SingleSheetDataEntry $ws 6 "Name" "Address" "Contact"

# Row 7: entry with the Address column left blank
Add-SingleSheetDataEntry $ws 7 "Name" "" "Contact"
